$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.061724
$ws.Range("H2").Value = 0.185172
$ws.Range("I2").Value = 0.09652262708432048
$ws.Range("J2").Value = 0.09652262708432047
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 0.9725011027853334
$ws.Range("R2").Value = 8.752509925068001
$ws.Range("S2").Value = 0.03108566140127806
$ws.Range("T2").Value = 0.03108566140127806

$ws.Range("G3").Value = 0.061724
$ws.Range("H3").Value = 0.185172
$ws.Range("I3").Value = 0.09652262708432048
$ws.Range("J3").Value = 0.09652262708432047
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 1.663699272802667
$ws.Range("R3").Value = 14.973293455224
$ws.Range("S3").Value = 0.0531795718480662
$ws.Range("T3").Value = 0.05317957184806617

$ws.Range("G4").Value = 0.061724
$ws.Range("H4").Value = 0.185172
$ws.Range("I4").Value = 0.09652262708432048
$ws.Range("J4").Value = 0.09652262708432047
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 0.3834671190653333
$ws.Range("R4").Value = 3.451204071588
$ws.Range("S4").Value = 0.01225739383497623
$ws.Range("T4").Value = 0.01225739383497623

$ws.Range("I5").Value = 0.8735221647273214
$ws.Range("J5").Value = 0.8735221647273215
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 8.801058302760888
$ws.Range("R5").Value = 79.209524724848
$ws.Range("S5").Value = 0.2813227847135126
$ws.Range("T5").Value = 0.2813227847135126

$ws.Range("I6").Value = 0.8735221647273214
$ws.Range("J6").Value = 0.8735221647273215
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.4812709322490146
$ws.Range("T6").Value = 0.4812709322490145

$ws.Range("I7").Value = 0.8735221647273214
$ws.Range("J7").Value = 0.8735221647273215
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.1109284477647943
$ws.Range("T7").Value = 0.1109284477647943

$ws.Range("I8").Value = 0.02995520818835809
$ws.Range("J8").Value = 0.02995520818835809
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 0.3018097815747778
$ws.Range("R8").Value = 2.716288034173
$ws.Range("S8").Value = 0.009647245284099357
$ws.Range("T8").Value = 0.009647245284099356

$ws.Range("I9").Value = 0.02995520818835809
$ws.Range("J9").Value = 0.02995520818835809
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.01650395554075573
$ws.Range("T9").Value = 0.01650395554075572

$ws.Range("I10").Value = 0.02995520818835809
$ws.Range("J10").Value = 0.02995520818835809
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("S10").Value = 0.003804007363503008
$ws.Range("T10").Value = 0.003804007363503008

